$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").WrapText = $false
$ws.Range("D3").WrapText = $false
$ws.Range("D18").WrapText = $false
$ws.Range("D19").WrapText = $false
$ws.Range("D20").WrapText = $false
$ws.Range("D21").WrapText = $false
$ws.Range("D26").WrapText = $false
$ws.Range("D27:D36").WrapText = $false
